# Converts an "RRGGBB" hex string into the Int64 COLORREF value that the
# PowerPoint COM object model expects for RGBColor.RGB (0x00BBGGRR).
function RGBFromHex($hex) {
    $r = [Convert]::ToInt64($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt64($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt64($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# --- 1. Table on slide 5: switch the applied table style -------------------
$s5 = $p.Slides.Item(5)
$tbl = $s5.Shapes.Item(2).Table
$tbl.ApplyStyle("{F3ADC762-E7AE-40D8-A3BA-FF9DF00BA4F1}")

# --- 2. Presentation theme: recolor to the "Office Theme" palette ----------
# (dk1/lt1 stay black/white; dk2, lt2, the 6 accents, hlink and folHlink are
#  updated to match the Office default theme colors.)
$tcs = $p.Slides.Item(1).ThemeColorScheme

$tcs.Item(3).RGB  = RGBFromHex("44546A")  # dk2
$tcs.Item(4).RGB  = RGBFromHex("E7E6E6")  # lt2
$tcs.Item(5).RGB  = RGBFromHex("5B9BD5")  # accent1
$tcs.Item(6).RGB  = RGBFromHex("ED7D31")  # accent2
$tcs.Item(7).RGB  = RGBFromHex("A5A5A5")  # accent3
$tcs.Item(8).RGB  = RGBFromHex("FFC000")  # accent4
$tcs.Item(9).RGB  = RGBFromHex("4472C4")  # accent5
$tcs.Item(10).RGB = RGBFromHex("70AD47")  # accent6
$tcs.Item(11).RGB = RGBFromHex("0563C1")  # hlink
$tcs.Item(12).RGB = RGBFromHex("954F72")  # folHlink
